# Apply the 2022-06-15 data update to the Fonds de solidarite worksheet.
# For a set of rows, update column C (nombre_aides) and column E (montant_total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column indices: C = 3, E = 5
$updates = @(
    @{ Row = 3;   C = 249335;  E = 1036488011 },
    @{ Row = 6;   C = 20796;   E = 360718555 },
    @{ Row = 7;   C = 7018;    E = 291118761 },
    @{ Row = 53;  C = 141687;  E = 590077232 },
    @{ Row = 82;  C = 8453;    E = 124867404 },
    @{ Row = 92;  C = 409316;  E = 1597549233 },
    @{ Row = 93;  C = 209672;  E = 1310272697 },
    @{ Row = 94;  C = 94241;   E = 919287786 },
    @{ Row = 95;  C = 50809;   E = 934729383 },
    @{ Row = 96;  C = 17327;   E = 797769865 },
    @{ Row = 166; C = 35931;   E = 210610514 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
